# Fix exploration time & agent_step_time formula (#13):
# recalculated Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M) and Std_Experiment_Time (N) for each data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.949490429999999
$ws.Range("H2").Value = 512.77234832
$ws.Range("M2").Value = 1.05778822273519
$ws.Range("N2").Value = 92.91633951288924
$ws.Range("G3").Value = 8.108459160000001
$ws.Range("H3").Value = 708.0508983
$ws.Range("M3").Value = 0.7822569855540713
$ws.Range("N3").Value = 118.8119675302528
$ws.Range("G4").Value = 3.89035115
$ws.Range("H4").Value = 128.60162033
$ws.Range("M4").Value = 0.7040189259567408
$ws.Range("N4").Value = 37.86378018965771
$ws.Range("G5").Value = 3.50247937
$ws.Range("H5").Value = 171.25943988
$ws.Range("M5").Value = 0.4861583988630538
$ws.Range("N5").Value = 41.27491655830466
$ws.Range("G6").Value = 1.47205695
$ws.Range("H6").Value = 29.0620712
$ws.Range("M6").Value = 0.3741776653142428
$ws.Range("N6").Value = 11.3149706709528
$ws.Range("G7").Value = 1.31934387
$ws.Range("H7").Value = 36.53761601999999
$ws.Range("M7").Value = 0.2532566510775939
$ws.Range("N7").Value = 11.91472697205999
$ws.Range("G8").Value = 0.7416628399999999
$ws.Range("H8").Value = 10.34928998
$ws.Range("M8").Value = 0.2489281767110648
$ws.Range("N8").Value = 4.81067637557166
$ws.Range("G9").Value = 0.6810248899999999
$ws.Range("H9").Value = 13.60453525
$ws.Range("M9").Value = 0.189282631053159
$ws.Range("N9").Value = 6.493293347756063
$ws.Range("G10").Value = 0.41474409
$ws.Range("H10").Value = 4.38585765
$ws.Range("M10").Value = 0.1744700250455972
$ws.Range("N10").Value = 2.657549453933808
$ws.Range("G11").Value = 0.3976717
$ws.Range("H11").Value = 6.285023280000001
$ws.Range("M11").Value = 0.1362788077418846
$ws.Range("N11").Value = 3.676985495048388
$ws.Range("G12").Value = 0.25240488
$ws.Range("H12").Value = 2.08240106
$ws.Range("M12").Value = 0.1167621609049987
$ws.Range("N12").Value = 1.417590722499761
$ws.Range("G13").Value = 0.22584488
$ws.Range("H13").Value = 2.91814978
$ws.Range("M13").Value = 0.08359638611780749
$ws.Range("N13").Value = 1.936955307406405
